$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last refreshed" timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 5 de Julio de 2020 a las 18:01"

# --- Refreshed COVID country statistics ---
# Columns: A Pais, B Casos totales, C Nuevos casos, D Casos activos,
#          E Recuperados, F Casos criticos, G Muertes hoy, H Muertes
$rowsData = @{
    4   = @("Estados Unidos", 2953014, 17244, 1260739, 1559893, 0, 64, 132382)
    # India overtakes Rusia for the #3 spot -> rows 6/7 swap
    6   = @("India", 687760, 13856, 418945, 249247, 0, 289, 19568)
    7   = @("Rusia", 681251, 6736, 450750, 220340, 0, 134, 10161)
    10  = @("Chile", 295532, 3685, 261032, 28192, 0, 116, 6308)
    11  = @("Reino Unido", 285416, 516, 0, 0, 0, 22, 44220)
    13  = @("Italia", 241611, 192, 192108, 14642, 0, 7, 34861)
    # Azerbaiyan overtakes Ghana -> rows 57/58 swap
    57  = @("Azerbaiyan", 20324, 523, 11742, 8332, 0, 9, 250)
    58  = @("Ghana", 20085, 697, 14870, 5093, 0, 5, 122)
    69  = @("Chequia", 12469, 29, 7864, 4257, 0, 0, 348)
    95  = @("Luxemburgo", 4522, 46, 4016, 396, 0, 0, 110)
    99  = @("Grecia", 3519, 8, 1374, 1953, 0, 0, 192)
    131 = @("Jordania", 1164, 14, 902, 252, 0, 0, 10)
    # Fiyi / Dominica tied on totals -> rows 205/206 swap (values unchanged)
    205 = @("Fiyi", 18, 0, 18, 0, 0, 0, 0)
    206 = @("Dominica", 18, 0, 18, 0, 0, 0, 0)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
    }
}
